$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Updated Wilke reference test results (AMS [N] measurements in column C).
# The D column formulas (=Cn/C5) recalculate automatically.
$ws.Range("C2").Value = 111.4
$ws.Range("C3").Value = 305.6
$ws.Range("C4").Value = 472.2
$ws.Range("C5").Value = 582.7
$ws.Range("C6").Value = 1216.8
$ws.Range("C7").Value = 1247.2
$ws.Range("C8").Value = 2528.4
$ws.Range("C9").Value = 2159.2

# Move/record the active selection as it was when the workbook was last saved.
$ws.Range("J20").Select()
